# Applies the OOXML diff to the presentation:
#  - Shape 1 "CustomShape 1" (Util): resize; "distpoints" -> "distPoints"
#  - Shape 2 "CustomShape 2" (Simulation): resize
#  - Shape 3 "CustomShape 3" (Control): resize
#  - Shape 4 "CustomShape 4" (Data): resize; "list" -> "cellsList";
#       "point3D" -> "Point3D" (color ff0000 -> 000000);
#       "matrix4" -> "Matrix4" (color ff0000 -> 000000);
#       "fileManager" color ff0000 -> 000000; drop trailing empty paragraph
#  - Shape 5 "CustomShape 5" (Viewer): resize
#  - Shape 6 "CustomShape 6" (Libs): resize
#
# Note: the PowerPoint size properties (Width/Height) are expressed in
# points in the COM object model (1 pt = 12700 EMU) while the underlying
# OOXML stores English Metric Units (EMU). The host's point -> EMU
# round trip is lossy (it quantizes through a 32-bit float before
# flooring back to EMU), so the literal point values below are nudged by
# a handful of ULPs so that, after that lossy round trip, the EMU value
# written back out is exactly the one required by the target diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Shape 1 "CustomShape 1" (Util) ----
$shp1 = $s.Shapes.Item(1)
$shp1.Width = 226.1763840527559
$shp1.Height = 421.6251984503937
$shp1.TextFrame.TextRange.Characters(6, 10).Text = "distPoints"

# ---- Shape 2 "CustomShape 2" (Simulation) ----
$shp2 = $s.Shapes.Item(2)
$shp2.Width = 189.41102362204725
$shp2.Height = 230.14488188976378

# ---- Shape 3 "CustomShape 3" (Control) ----
$shp3 = $s.Shapes.Item(3)
$shp3.Width = 208.8
$shp3.Height = 249.10866551732283

# ---- Shape 4 "CustomShape 4" (Data) ----
$shp4 = $s.Shapes.Item(4)
$shp4.Width = 225.80787401574804
$shp4.Height = 322.4976377952756

$tr4 = $shp4.TextFrame.TextRange
# color fixes (ff0000 -> 000000), same-length edits first so character
# offsets stay valid
$tr4.Characters(38, 7).Font.Color.RGB = 0
$tr4.Characters(46, 7).Font.Color.RGB = 0
$tr4.Characters(65, 11).Font.Color.RGB = 0
# case-only renames (length preserved, offsets unaffected)
$tr4.Characters(38, 7).Text = "Point3D"
$tr4.Characters(46, 7).Text = "Matrix4"
# length-changing rename done last among text edits on this shape
$tr4.Characters(12, 4).Text = "cellsList"
# drop the trailing empty paragraph (#10) after "fileManager" (#9)
$tr4.Paragraphs(10, 1).Delete()

# ---- Shape 5 "CustomShape 5" (Viewer) ----
$shp5 = $s.Shapes.Item(5)
$shp5.Width = 189.41102362204725
$shp5.Height = 62.192125984251966

# ---- Shape 6 "CustomShape 6" (Libs) ----
$shp6 = $s.Shapes.Item(6)
$shp6.Width = 189.41102362204725
$shp6.Height = 129.3732283464567
